$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Smartphone Multi G Max 2"
$ws.Range("C2").Value = "971,91"

$ws.Range("B3").Value = "Smartphone Multi F"
$ws.Range("C3").Value = "533,61"

$ws.Range("B4").Value = "Smartphone Galaxy A04e"
$ws.Range("C4").Value = "854,91"

$ws.Range("B5").Value = "Smartphone Multi G Max 2"
$ws.Range("C5").Value = "809,91"

$ws.Range("B7").Value = "Smartphone Multi E P9128"
$ws.Range("C7").Value = "539,10"

$ws.Range("B8").Value = "Smartphone Multi G 2"
$ws.Range("C8").Value = "728,91"

$ws.Range("B9").Value = "Smartphone Multi G"
$ws.Range("C9").Value = "614,61"

$ws.Range("B10").Value = "Smartphone Multilaser G 32gb 5mp 5.5PT P9132 Multi CX 1 UN"
$ws.Range("C10").Value = "614,61"
